$wb = $excel.ActiveWorkbook

# --- Input sheet: just move the active selection to G2 ---
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Activate()
$wsInput.Range("G2").Select()

# --- Output sheet: add a new value in D2 and move selection there ---
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Activate()
$wsOutput.Range("D2").Value = "y"
$wsOutput.Range("D2").Select()

# --- Joins sheet: remove the REVISIONS/DID join row (row 2), which moves
#     the DOCUMENTS/DORIGINALNAME row up, then flip the Join flag to "n" ---
$wsJoins = $wb.Worksheets.Item("Joins")
$wsJoins.Activate()
$wsJoins.Rows.Item(2).Delete()
$wsJoins.Range("E2").Value = "n"
